$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-48: update Price (D) and Volume(1h) (E) values ---
# D-column cells hold text-formatted numbers (e.g. "321.82", "43.475.09").
# A leading apostrophe forces Excel to store the literal text instead of
# auto-converting to a numeric value (which would also strip formatting such
# as trailing zeros). Resetting the Style back to "Normal" afterwards clears
# the quote-prefix marker Excel applies, matching the original cell styling.
$ws.Range("D2").Value = "'43.475.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "'2.242.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'321.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "'100.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "'36.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "'0.0829"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "'7.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").Value = "'2.581.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "'0.854"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'14.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("D17").Value = "'2.239.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'43.379.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").Value = "'13.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.61%  "
$ws.Range("D20").Value = "'0.0₃0983"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").Value = "'6.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'65.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "'236.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("D29").Value = "'36.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("D30").Value = "'6.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").Value = "'159.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").Value = "'20.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "'0.0850"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D34").Value = "'2.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").Value = "'3.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").Value = "'0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.42%  "
$ws.Range("D37").Value = "'1.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "'4.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("D41").Value = "'15.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.43%  "
$ws.Range("D42").Value = "'0.0316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "'1.794.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("D46").Value = "'82.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.29%  "
$ws.Range("D47").Value = "'74.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("E48").Value = "  -3.40%  "

# --- Rows 49-51: coins reordered (Stacks, MultiversX, Aave -> MultiversX, Aave, Stacks) ---
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'58.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'103.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.55%  "
